$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Re-arrange the header row (row 1): insert "Date of Registration
#    (RAC)", "Category" and "Gender" columns, shift Degree /
#    Awarded-Submitted / Thesis Title accordingly. Columns J:L are
#    brand new, so borrow A1's header style before typing into them.
# ---------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("J1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A1").Value = "Scholar Name"
$ws.Range("B1").Value = "Department Name"
$ws.Range("C1").Value = "Guide Name"
$ws.Range("D1").Value = "Thesis Title"
$ws.Range("E1").Value = "Date of Registration (RAC)"
$ws.Range("I1").Value = "Gender"
$ws.Range("G1").Value = "Degree"
$ws.Range("H1").Value = "Awarded / Submitted"
$ws.Range("F1").Value = "Category"
$ws.Range("J1").Value = "Year of Scholar Registration"
$ws.Range("K1").Value = "Year of Award"
$ws.Range("L1").Value = "Choose Year"

# ---------------------------------------------------------------
# 2. Data row (row 2). Give E2:G2 the same number-format style (the
#    author applied a date format across the block before typing the
#    Category / Degree values into F2/G2).
# ---------------------------------------------------------------
$ws.Range("A2").Value = "aa"
$ws.Range("B2").Value = "school"
$ws.Range("C2").Value = "aaaa"
$ws.Range("D2").Value = "aaaa"

$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2:G2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$regDate = Get-Date -Year 1999 -Month 11 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("E2").Value = $regDate
$ws.Range("I2").Value = "Male"
$ws.Range("G2").Value = "Ph.D."
$ws.Range("F2").Value = "OBC"
$ws.Range("H2").Value = "Awarded"
$ws.Range("J2").Value = 2020
$ws.Range("K2").Value = 2020
$ws.Range("L2").Value = "2019-20"

# Row 3 used to carry the "Choose Year" sample in column I; it now
# belongs under column L, so clear the stale I3 cell.
$ws.Range("I3").ClearContents() | Out-Null
$ws.Range("L3").Value = "2022-23"

# ---------------------------------------------------------------
# 3. Column widths / layout. Column D is inserted before the old
#    "Degree" column, so columns 3 & 4 share the old column-3 width,
#    and several columns further right get new widths.
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.7369791666667
$ws.Columns.Item(2).ColumnWidth = 22.1666666666667
$ws.Columns.Item(3).ColumnWidth = 16.3072916666667
$ws.Columns.Item(4).ColumnWidth = 16.3072916666667
$ws.Columns.Item(5).ColumnWidth = 28.8776041666667
$ws.Columns.Item(6).ColumnWidth = 13.5924479166667
$ws.Columns.Item(7).ColumnWidth = 11.7369791666667
$ws.Columns.Item(8).ColumnWidth = 22.3072916666667
$ws.Columns.Item(9).ColumnWidth = 9.59244791666667
$ws.Columns.Item(10).ColumnWidth = 28.3072916666667
$ws.Columns.Item(11).ColumnWidth = 16.7369791666667
$ws.Columns.Item(12).ColumnWidth = 13.0221354166667

# ---------------------------------------------------------------
# 4. Data validations - drop the old four and rebuild the seven that
#    the new layout needs, in the same order the author added them.
# ---------------------------------------------------------------
$ws.Range("E2:E1002").Validation.Delete() | Out-Null
$ws.Range("E1").Validation.Delete() | Out-Null
$ws.Range("G1:H1002").Validation.Delete() | Out-Null
$ws.Range("I2:I1048576").Validation.Delete() | Out-Null

$r = $ws.Range("J1:K1002")
$r.Validation.Add(2, 1, 5, 0, 0) | Out-Null
$r.Validation.ShowInput = $false

$r = $ws.Range("L2:L1048576")
$r.Validation.Add(3, 1, 1, '"1994-95,1995-96,1996-97,1997-98,1998-99,1999-20,2000-01,2001-02,2002-03,2003-04,2004-05,2005-06,2006-07,2007-08,2008-09,2009-10,2010-11,2011-12,2012-13,2013-14,2014-15,2015-16,2016-17,2017-18,2018-19,2019-20,2020-21,2021-22,2022-23,"') | Out-Null

$r = $ws.Range("E2:E1048576")
$r.Validation.Add(4, 1, 7, 35690) | Out-Null

$r = $ws.Range("I2:I1048576")
$r.Validation.Add(3, 1, 1, '"Male, Female, Other"') | Out-Null

$r = $ws.Range("G2:G1048576")
$r.Validation.Add(3, 1, 7, '"Ph.D., M.Phil, PG Dissertation"') | Out-Null

$r = $ws.Range("F2:F1048576")
$r.Validation.Add(3, 1, 7, '"Open, SC, ST,VJ / NT(A),VJ / NT(B),VJ / NT( C),VJ / NT(D),SBC,OBC,PH,EWS,Minority,NRI,Foreign National"') | Out-Null

$r = $ws.Range("H1:H1048576")
$r.Validation.Add(3, 1, 3, '"Awarded,Submitted,Ongoing"') | Out-Null
$r.Validation.ShowInput = $false

# ---------------------------------------------------------------
# 5. Selection / view - the author scrolled back to A1 and left the
#    cursor on H2 (Awarded / Submitted).
# ---------------------------------------------------------------
$ws.Range("H2").Select() | Out-Null
